$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (new most-recent quarters)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formatting from the (now-shifted) old column D/E -- now F/G -- onto the new D:E columns
$ws.Range("F5:G102").Copy()
$ws.Range("D5").PasteSpecial(-4122)

# Populate the two new columns with the latest two quarters of data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 1045600
$ws.Range("E8").Value = 1076500
$ws.Range("D9").Value = 87200
$ws.Range("E9").Value = 66800
$ws.Range("D10").Value = 958400
$ws.Range("E10").Value = 1009700
$ws.Range("D12").Value = 900
$ws.Range("E12").Value = 700
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 143400
$ws.Range("E14").Value = 222300
$ws.Range("D15").Value = 263000
$ws.Range("E15").Value = 243200
$ws.Range("D17").Value = 1092300
$ws.Range("E17").Value = 1071700
$ws.Range("D18").Value = -46700
$ws.Range("E18").Value = 4800
$ws.Range("D20").Value = 12500
$ws.Range("E20").Value = 10700
$ws.Range("D21").Value = 229500
$ws.Range("E21").Value = 259400
$ws.Range("D22").Value = 78400
$ws.Range("E22").Value = 74500
$ws.Range("D23").Value = -112600
$ws.Range("E23").Value = -59000
$ws.Range("D24").Value = -131400
$ws.Range("E24").Value = 19000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 18700
$ws.Range("E26").Value = -78000
$ws.Range("D27").Value = -121500
$ws.Range("E27").Value = -154400
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -12500
$ws.Range("E32").Value = -10700
$ws.Range("D33").Value = -121500
$ws.Range("E33").Value = -154400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -121500
$ws.Range("E35").Value = -154400
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 0
$ws.Range("E41").Value = 0
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 525900
$ws.Range("E43").Value = 400600
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 280700
$ws.Range("E45").Value = 506000
$ws.Range("D46").Value = 806600
$ws.Range("E46").Value = 906600
$ws.Range("D47").Value = 433600
$ws.Range("E47").Value = 392900
$ws.Range("D48").Value = 13869900
$ws.Range("E48").Value = 13784400
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 409300
$ws.Range("E52").Value = 718600
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 15519500
$ws.Range("E54").Value = 15802500
$ws.Range("D57").Value = 66300
$ws.Range("E57").Value = 91900
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 787300
$ws.Range("E59").Value = 721900
$ws.Range("D60").Value = 853500
$ws.Range("E60").Value = 813900
$ws.Range("D61").Value = 5461700
$ws.Range("E61").Value = 5487000
$ws.Range("D62").Value = 716800
$ws.Range("E62").Value = 830500
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 7853700
$ws.Range("E66").Value = 7888900
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 1177500
$ws.Range("E72").Value = 1299100
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 7665800
$ws.Range("E76").Value = 7913600
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -121500
$ws.Range("E81").Value = -154400
$ws.Range("D83").Value = 263700
$ws.Range("E83").Value = 243900
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 821600
$ws.Range("E89").Value = 421500
$ws.Range("D91").Value = -2100
$ws.Range("E91").Value = -1900
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -580400
$ws.Range("E94").Value = -598000
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -241200
$ws.Range("E100").Value = 125900
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 0
$ws.Range("E102").Value = -50600
